$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the "Förändrad" date column (C) for rows 2-7 from 45183 (2023-09-14)
# to 45184 (2023-09-15), matching the automatic daily update in the diff.
foreach ($row in 2..7) {
    $ws.Range("C$row").Value = 45184
}
